# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 7 de Julio de 2020 a las 17:30"

# Row 4 - Estados Unidos: refreshed totals
$ws.Range("B4").Value = 3050476
$ws.Range("C4").Value = 9643
$ws.Range("D4").Value = 1326669
$ws.Range("E4").Value = 1590682
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 146
$ws.Range("H4").Value = 133125

# Chile overtakes España in ranking -> rows 9 & 10 swap positions
# Row 9 becomes Chile (with refreshed data)
$ws.Range("A9").Value = "Chile"
$ws.Range("B9").Value = 301019
$ws.Range("C9").Value = 2462
$ws.Range("D9").Value = 264371
$ws.Range("E9").Value = 30214
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 50
$ws.Range("H9").Value = 6434

# Row 10 becomes España (previous Chile row's slot)
$ws.Range("A10").Value = "España"
$ws.Range("B10").Value = 298869
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 0
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 28388

# Row 28 - Suecia: refreshed totals
$ws.Range("B28").Value = 73344
$ws.Range("C28").Value = 57
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 11
$ws.Range("H28").Value = 5447

# Row 61 - Moldavia: refreshed totals
$ws.Range("E61").Value = 6062
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 11
$ws.Range("H61").Value = 603

# Row 146 - Montenegro: refreshed totals
$ws.Range("E146").Value = 510
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 2
$ws.Range("H146").Value = 16

# Row 154 - Surinam: refreshed totals
$ws.Range("E154").Value = 284
$ws.Range("F154").Value = 0
$ws.Range("G154").Value = 1
$ws.Range("H154").Value = 15

# Groenlandia overtakes Islas Malvinas in ranking -> rows 209 & 210 swap labels
# (underlying totals for both countries are identical, so only the names swap)
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"
